$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    # Prefix with an apostrophe so Excel stores the value as plain text
    # instead of auto-converting date-like / numeric-looking strings, then
    # reset the style back to Normal so no quotePrefix / number-format
    # style ends up applied to the cell (matching the original formatting).
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

# Row 2 (Alcoa)
Set-TextValue $ws.Range("A2") "2025/12/05"
Set-TextValue $ws.Range("B2") "5.39"

# Row 8 (Rio Tinto)
Set-TextValue $ws.Range("A8") "2025/12/05"
Set-TextValue $ws.Range("B8") "7.91"

# Row 14 (Norsk Hydro)
Set-TextValue $ws.Range("A14") "2025/12/05"
Set-TextValue $ws.Range("B14") "2.88"

# Row 20 (Reliance)
Set-TextValue $ws.Range("A20") "2025/12/05"
Set-TextValue $ws.Range("B20") "12.73"

# Row 26 (Kaiser)
Set-TextValue $ws.Range("A26") "2025/12/05"
Set-TextValue $ws.Range("B26") "10.60"

# Row 32 (Ryerson)
Set-TextValue $ws.Range("A32") "2025/12/05"
Set-TextValue $ws.Range("B32") "26.69"

# Row 38 (Alro Steel) - date only, EBITDA unchanged
Set-TextValue $ws.Range("A38") "2025/12/05"

# Row 44 (Ultra)
Set-TextValue $ws.Range("A44") "2025/12/05"
Set-TextValue $ws.Range("B44") "11.30"

# Row 50 (Benchmark)
Set-TextValue $ws.Range("A50") "2025/12/05"
Set-TextValue $ws.Range("B50") "12.08"

# Row 56 (Celestica)
Set-TextValue $ws.Range("A56") "2025/12/05"
Set-TextValue $ws.Range("B56") "33.78"

# Row 62 (Jabil)
Set-TextValue $ws.Range("A62") "2025/12/05"
Set-TextValue $ws.Range("B62") "11.86"

# Row 68 (Flex)
Set-TextValue $ws.Range("A68") "2025/12/05"
Set-TextValue $ws.Range("B68") "12.76"

# Row 74 (MKS) - date only, EBITDA unchanged
Set-TextValue $ws.Range("A74") "2025/12/05"
